$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.768.59'
$ws.Range('E2').Value = '  +1.27%  '

$ws.Range('D3').Value = '3.155.01'
$ws.Range('E3').Value = '  +0.86%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = '532.61'
$ws.Range('E5').Value = '  +0.27%  '

$ws.Range('D6').Value = '140.43'
$ws.Range('E6').Value = '  +1.26%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E9').Value = '  +0.54%  '

$ws.Range('D10').Value = '0.434'
$ws.Range('E10').Value = '  +6.20%  '

$ws.Range('E11').Value = '  +3.52%  '

$ws.Range('E12').Value = '  +2.69%  '

$ws.Range('D13').Value = '3.700.06'
$ws.Range('E13').Value = '  +0.98%  '

$ws.Range('D14').Value = '26.23'
$ws.Range('E14').Value = '  +2.88%  '

$ws.Range('D15').Value = '0.0000173'
$ws.Range('E15').Value = '  +5.74%  '

$ws.Range('D16').Value = '58.818.27'
$ws.Range('E16').Value = '  +1.34%  '

$ws.Range('E17').Value = '  +4.57%  '

$ws.Range('D18').Value = '3.150.45'
$ws.Range('E18').Value = '  +1.35%  '

$ws.Range('D19').Value = '13.02'
$ws.Range('E19').Value = '  +2.78%  '

$ws.Range('D20').Value = '8.20'
$ws.Range('E20').Value = '  +1.29%  '

$ws.Range('D21').Value = '372.48'
$ws.Range('E21').Value = '  +5.44%  '

$ws.Range('E22').Value = '  +2.02%  '

$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  +0.02%  '

$ws.Range('D24').Value = '70.08'
$ws.Range('E24').Value = '  +1.95%  '

$ws.Range('D25').Value = '0.519'
$ws.Range('E25').Value = '  +3.07%  '

$ws.Range('D26').Value = '0.167'
$ws.Range('E26').Value = '  +0.23%  '

$ws.Range('E27').Value = '  -0.11%  '

$ws.Range('D28').Value = '8.20'
$ws.Range('E28').Value = '  +12.77%  '

$ws.Range('D29').Value = '0.0₃0866'
$ws.Range('E29').Value = '  -2.18%  '

$ws.Range('E30').Value = '  +1.42%  '

$ws.Range('E31').Value = '  -0.10%  '

$ws.Range('D32').Value = '22.16'
$ws.Range('E32').Value = '  +4.01%  '

$ws.Range('D33').Value = '5.20'
$ws.Range('E33').Value = '  +4.07%  '

$ws.Range('E34').Value = '  +1.02%  '

$ws.Range('D35').Value = '159.09'
$ws.Range('E35').Value = '  +0.20%  '

$ws.Range('D36').Value = '6.28'
$ws.Range('E36').Value = '  +3.51%  '

$ws.Range('E37').Value = '  +7.03%  '

$ws.Range('D38').Value = '25.24'
$ws.Range('E38').Value = '  -2.99%  '

$ws.Range('E39').Value = '  -0.40%  '

$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.0685'
$ws.Range('E40').Value = '  +2.14%  '

$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.629.90'
$ws.Range('E41').Value = '  +9.94%  '

$ws.Range('E42').Value = '  +6.13%  '

$ws.Range('E43').Value = '  +7.88%  '

$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '38.80'
$ws.Range('E44').Value = '  +3.18%  '

$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '0.712'
$ws.Range('E45').Value = '  +1.90%  '

$ws.Range('E46').Value = '  +0.08%  '

$ws.Range('D47').Value = '3.197.14'
$ws.Range('E47').Value = '  +1.00%  '

$ws.Range('E48').Value = '  +14.42%  '

$ws.Range('E49').Value = '  +0.70%  '

$ws.Range('D50').Value = '6.20'
$ws.Range('E50').Value = '  +2.85%  '

$ws.Range('D51').Value = '20.27'
$ws.Range('E51').Value = '  +2.14%  '
